# Applies two changes described by the commit diff:
#   1. Re-point the single table on the deck (slide 16) at a different
#      table style (tableStyleId GUID swap).
#   2. Re-colour the presentation's theme (theme1.xml, the slide master's
#      theme - "Integral") so its 12 theme colours become the stock
#      "Office Theme" palette.
#
# (theme2.xml, the separate Notes Master theme, is swapped the opposite
# way in the source diff, but the Notes Master's theme is not reachable
# as a distinct object through the PowerPoint COM object model - it
# resolves back to the same Theme as the slide master - so it cannot be
# targeted independently here.)

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$newTableStyleId = "{3FDB56A7-4EA2-4889-A703-DB4EE8D07AB2}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colours (theme1.xml: "Integral" -> "Office Theme") ----
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# index : name     : target RGB (Office Theme palette)
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
